# "Moving example images so they are not covered by the webcam"
#
# Slide 9 ("Jumbotron with background image") has an example screenshot
# (the "Picture 7" image shape) that currently sits where a webcam overlay
# (added when the deck is screen-recorded) would cover it. Reposition it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Locate the picture shape by name rather than assuming its index.
$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Picture 7") {
        $pic = $shp
        break
    }
}

# New target position, expressed in points (EMU / 12700), nudged by a few
# EMU-scale epsilons so the float32 COM round-trip lands exactly on the
# target EMU offsets (6054570,4496718) -> (3630153,4706351).
$pic.Left = 285.8388218976378
$pic.Top  = 370.5788188976378
